$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 corresponds to file "metrics_sim_with_priors.json"
$ws.Range("C3").Value = 0.8828828828828829

$ws.Range("H3").Value = 0.763432446531038
$ws.Range("I3").Value = 0.02934454899520402
$ws.Range("J3").Value = 0.7837837837837838
$ws.Range("K3").Value = 277.7567567567568

$ws.Range("Q3").Value = 4
$ws.Range("R3").Value = 12
$ws.Range("S3").Value = 32
$ws.Range("T3").Value = 195
$ws.Range("U3").Value = 672
$ws.Range("V3").Value = 7553
$ws.Range("W3").Value = 7545
$ws.Range("X3").Value = 7525
$ws.Range("Y3").Value = 7362
$ws.Range("Z3").Value = 6885

$ws.Range("AF3").Value = 0.999471
$ws.Range("AG3").Value = 0.998412
$ws.Range("AH3").Value = 0.995766
$ws.Range("AI3").Value = 0.974196
$ws.Range("AJ3").Value = 0.911076
